$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text update (shared string reused across sheets) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet: handback info for this run ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("I2").Value = "ad524f1c-9d01-4152-ac78-f7db533094fb.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1b489718c3c7f5a1b16f9a9767e72eb07165d24/e2e/ad524f1c-9d01-4152-ac78-f7db533094fb.md", "", "", "ad524f1c-9d01-4152-ac78-f7db533094fb.md")
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("J2").Value = "ad524f1c-9d01-4152-ac78-f7db533094fb.703042b5688bf47aa2cf97d1f0a04d1ec3831e5d.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-26 06:57:34"
$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# --- de-de sheet: handback info for this run ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("I2").Value = "ad524f1c-9d01-4152-ac78-f7db533094fb.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1b489718c3c7f5a1b16f9a9767e72eb07165d24/e2e/ad524f1c-9d01-4152-ac78-f7db533094fb.md", "", "", "ad524f1c-9d01-4152-ac78-f7db533094fb.md")
$wsDe.Range("I2").Style = "HyperLink"
$wsDe.Range("J2").Value = "ad524f1c-9d01-4152-ac78-f7db533094fb.703042b5688bf47aa2cf97d1f0a04d1ec3831e5d.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-26 06:57:40"
$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
